# Sanity Semilla 4 - add "Semilla 3" sheet, re-point selections.
#
# The original workbook has a single sheet "Semilla 4". This change:
#   1. Duplicates it as a new sheet "Semilla 3" placed after "Semilla 4".
#   2. On the new sheet, re-targets the EPOS/CRM/CONFIRMADOR hyperlinks + values
#      in row 2, and the DB host / hostname values in rows 4-7, to a different
#      environment (kept as the active/selected sheet, cell E9 selected).
#   3. Leaves "Semilla 4" itself selected at cell B20 (no longer the active tab).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Duplicate "Semilla 4" right after itself; this clones values, styles,
#    hyperlinks, column widths, etc. exactly.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Semilla 3"

# 2. Re-create the 5 hyperlinks on the new sheet so the 3 that change
#    (EPOS/CRM/CONFIRMADOR, in A2/B2/C2) point at the new environment while
#    GATEWAYCBS/GATEWAY MG (D2/E2) keep their original targets. Re-adding
#    (rather than editing in place) avoids leaving stale relationships
#    around, and doing all 5 keeps the relationship ids contiguous.
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("D2"), "http://10.65.45.12:9001/gatewaycbs/BcServicesInt")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "http://10.65.45.12:9001/gatewaymgint/GatewayMGWSInt")
$ws2.Hyperlinks.Add($ws2.Range("A2"), "http://10.69.60.106:8180/tigo-pos-web/index.jsp")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "http://10.69.60.107:8080/CRMPortal/auth/portal/default/Venta")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "http://10.69.60.106:8180/tigo-pos-web/wap/windex.wml")

# Adding hyperlinks stamps its own "hyperlink" cell style; restore the
# original formatting (inherited from "Semilla 4" row 2) over it.
$ws1.Range("A2:H2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Now write the new, environment-specific text for row 2 + the DB rows.
$ws2.Range("C2").Value = "http://10.69.60.106:8180/tigo-pos-web/wap/windex.wml"
$ws2.Range("B2").Value = "http://10.69.60.107:8080/CRMPortal/auth/portal/default/Venta"
$ws2.Range("A2").Value = "http://10.69.60.106:8180/tigo-pos-web/index.jsp"

$ws2.Range("A4").Value = "10.69.60.103"
$ws2.Range("A5").Value = "10.69.60.102"
$ws2.Range("A6").Value = "10.69.60.102"
$ws2.Range("A7").Value = "10.65.32.76"
$ws2.Range("B7").Value = "SIEBEL02"

# 3. Selection bookkeeping: "Semilla 4" keeps cell B20 selected (and is no
#    longer the active tab); the new "Semilla 3" ends up active with E9
#    selected. Select on ws1 first so the final Select() on ws2 is what
#    leaves it as the active sheet/tab.
$ws1.Range("B20").Select()
$ws2.Range("E9").Select()
